# Add a "getCOTaskTable"-related pair of rows (tgtID / bumpID) to the
# trial-table label sheet.
#
# New rows appended to the bottom of the table:
#   19: tgtID   | integer indicating target code (1-16)
#   20: bumpID  | integer indicating bump code (1-16)
#
# The order in which the cell values are assigned below controls the order
# new entries are appended to the shared string table, matching the
# target workbook (bumpID, "integer indicating target code (1-16)",
# "integer indicating bump code (1-16)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "bumpID"
$ws.Range("B19").Value = "integer indicating target code (1-16)"
$ws.Range("B20").Value = "integer indicating bump code (1-16)"
$ws.Range("A19").Value = "tgtID"

# Move / record the active selection like Excel would after typing the
# data in and pressing enter past the last new row.
$ws.Range("B21").Select() | Out-Null
